$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").ClearFormats()

$text = @'
questions = [
    {
        "title": "Your organization wants to set up an Azure machine learning (ML) infrastructure that supports ML workflows for multiple teams in different regions with minimal manual intervention. Each team requires different types of virtual machines (VMs) for their workloads, and you want to ensure that the solution is cost effective.Which strategy should you use?",
        "ques_type": 2,
        "options": [
            "Configure an Azure Policy to limit the number of VMs that each team can use.",
            "Create separate Azure ML Compute clusters for each team and configure auto-scaling and auto-pause.",
            "Use a single Azure ML Compute cluster and adjust the VM size manually based on the workload.",
            "Set up Azure Databricks workspaces for each team."
        ],
        "score": "Create separate Azure ML Compute clusters for each team and configure auto-scaling and auto-pause."
    },
    {
        "title": "You are developing a complex machine learning (ML) workflow in Azure ML Designer. Your data science team is split over whether to use Principal Component Analysis (PCA) or t-Distributed Stochastic Neighbor Embedding (t-SNE) for dimensionality reduction, given that your high-dimensional dataset contains nonlinear structures.Which approach should you use?",
        "ques_type": 2,
        "options": [
            "Use PCA only.",
            "Use t-SNE only.",
            "Use PCA first, followed by t-SNE.",
            "Use t-SNE first, followed by PCA."
        ],
        "score": "Use PCA first, followed by t-SNE."
    },
    {
        "title": "You are developing a machine learning model for a global retail chain. The model will predict stock requirements based on real-time data from sensors located in the stores. Considering the network latency and connectivity issues, you want the predictions to be made at the store location in real time, without having to send the data back to a central location.Which Azure deployment option should you use?",
        "ques_type": 2,
        "options": [
            "Azure Kubernetes Service (AKS)",
            "Azure Container Instances (ACI)",
            "Azure Functions",
            "Azure IoT Edge"
        ],
        "score": "Azure IoT Edge"
    },
    {
        "title": "Your team is developing machine learning (ML) models using Azure ML, while your DevOps team uses Azure Pipelines for infrastructure as code (IaC), testing, and deployment. You need to develop a strategy to streamline the handoff between these teams while maintaining visibility into model performance and data lineage.Which approach should you use to ensure smooth operations?",
        "ques_type": 2,
        "options": [
            "Standardize on a single platform by having the data science team also use Azure Pipelines for model development.",
            "Use Azure Pipelines to deploy Azure ML models and track them with Azure ML model management capabilities.",
            "Have the data science team hand over the trained model to the DevOps team for deployment.",
            "Use Azure ML for both model development and deployment to minimize the use of Azure Pipelines."
        ],
        "score": "Use Azure Pipelines to deploy Azure ML models and track them with Azure ML model management capabilities."
    }
]
'@

$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit()
$ws.Range("A2").ClearContents()
